$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refit the "100 iterations" block (rows 2-11): updated D/E/F values ---
$ws.Range("D2").Value = 0.4091005477142459
$ws.Range("E2").Value = 0.4091005477142459

$ws.Range("D3").Value = 0.4681202610029285
$ws.Range("E3").Value = 0.4681202610029285

$ws.Range("D4").Value = 0.606438487890233
$ws.Range("E4").Value = 0.606438487890233

$ws.Range("D5").Value = 0.0000001121320758736525
$ws.Range("E5").Value = 0.0000001121320758736525

$ws.Range("D6").Value = 0.339944850596451
$ws.Range("E6").Value = 0.339944850596451

$ws.Range("D7").Value = 0.8400159708232472
$ws.Range("E7").Value = 0.1599840291767528

$ws.Range("D8").Value = 0.4225270629620901
$ws.Range("E8").Value = 0.5774729370379099

$ws.Range("D9").Value = 0.396905171725738
$ws.Range("E9").Value = 0.603094828274262

$ws.Range("D10").Value = 0.6036118607889198
$ws.Range("E10").Value = 0.3963881392110802

$ws.Range("D11").Value = 0.4496533255199429
$ws.Range("E11").Value = 0.5503466744800571
$ws.Range("F11").Value = 0.5769394040107727

# --- New "Label" column (H): 0 = Control patient, 1 = MDD patient ---
$ws.Range("H1").Value = "Label"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1").VerticalAlignment = -4160
$ws.Range("H1").Borders.LineStyle = 1

$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1

$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
